$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = "539b9bf3-0642-4945-a00e-31319c9cb371"
$ws.Range("B8").Value = "2025-09-28T00:00:00.000Z"
$ws.Range("C8").Value = "T1"
$ws.Range("D8").Value = "transchool"
$ws.Range("G8").Value = "Pascal"
$ws.Range("H8").Value = "lebon"
$ws.Range("I8").Value = "bangouraibrahima57@gmail.com"
$ws.Range("K8").Value = "Oui"
$ws.Range("L8").Value = "Oui"

# Row 9
$ws.Range("A9").Value = "3e7cfffc-744c-4772-9898-c2a2e7dcefea"
$ws.Range("B9").Value = "2025-09-28T00:00:00.000Z"
$ws.Range("C9").Value = "T2"
$ws.Range("D9").Value = "transchool"
$ws.Range("G9").Value = "Pascal"
$ws.Range("H9").Value = "lebon"
$ws.Range("I9").Value = "ptijjo@hotmail.com"
$ws.Range("K9").Value = "Oui"
$ws.Range("L9").Value = "Oui"
